## Generate Report for Handback
## - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
##   on the Overview sheet (zh-cn / de-de status columns) and on each of the
##   per-locale detail sheets.
## - The per-locale "Latest Handback DateTime" is refreshed to the time the
##   handback report was generated.
## - The stale "handback file is not the latest" error message is cleared
##   now that the handback is in sync, for both locales.
## - Column widths for the (now longer/shorter) Status / Error Detail
##   columns are refreshed to fit the new text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns (E, F) to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-21 08:54:39"
$wsZhCn.Range("P2").Value = ""

# Status column (C) needs to be wider, Error Detail (P) can shrink back down.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719540550566

# ---- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-21 08:54:45"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719540550566
